$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player names per the official draw ("sorteo oficial") results.
$ws.Range("B2").Value = "Armada"
$ws.Range("C2").Value = "Papu"
$ws.Range("G2").Value = "Fale"
$ws.Range("H2").Value = "Coquina"

$ws.Range("B3").Value = "Tony"
$ws.Range("C3").Value = "Kike"
$ws.Range("G3").Value = "Ruso"

$ws.Range("B4").Value = "Palop"
$ws.Range("C4").Value = "Kero"
$ws.Range("G4").Value = "Lope"
$ws.Range("H4").Value = "Puche"
